$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the PAT-coefficient cells in the top matrix (rows 2-8) ---
$ws.Range("C2").Value = 10
$ws.Range("D3").Value = 10
$ws.Range("E4").Value = 10
$ws.Range("F5").Value = 10
$ws.Range("G6").Value = 10
$ws.Range("H7").Value = 10
$ws.Range("H8").Value = 60
$ws.Range("I8").Value = 10

# --- Remove the old row 9 (facebook) entirely ---
$ws.Range("A9:I9").ClearContents()

# --- Rewrite the lower matrix block, shifted up by one row (old 15-22 -> new 14-21) ---
# Header row (was row 15, now row 14)
$ws.Range("B14").Value = "X"
$ws.Range("C14").Value = "A"
$ws.Range("D14").Value = "B"
$ws.Range("E14").Value = "C"
$ws.Range("F14").Value = "D"
$ws.Range("G14").Value = "E"
$ws.Range("H14").Value = "F"
$ws.Range("I14").Value = "Y"

# Row 15 (was row 16)
$ws.Range("A15").Value = "X"
$ws.Range("B15").Value = 0
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0

# Row 16 (was row 17)
$ws.Range("A16").Value = "A"
$ws.Range("B16").Value = 0
$ws.Range("C16").Value = 0
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0

# Row 17 (was row 18)
$ws.Range("A17").Value = "B"
$ws.Range("B17").Value = 0
$ws.Range("C17").Value = 0
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0

# Row 18 (was row 19) -- N19/O19 formulas stay put, untouched
$ws.Range("A18").Value = "C"
$ws.Range("B18").Value = 0
$ws.Range("C18").Value = 0
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0

# Row 19 (was row 20) -- N20/O20 formulas stay put, untouched
$ws.Range("A19").Value = "D"
$ws.Range("B19").Value = 0
$ws.Range("C19").Value = 0
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 1
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0

# Row 20 (was row 21)
$ws.Range("A20").Value = "E"
$ws.Range("B20").Value = 0
$ws.Range("C20").Value = 0
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = 1
$ws.Range("I20").Value = 0

# Row 21 (was row 22)
$ws.Range("A21").Value = "F"
$ws.Range("B21").Value = 0
$ws.Range("C21").Value = 0
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 0.86
$ws.Range("I21").Value = 0.14

# --- Drop the now-obsolete old rows 22 and 23 (facebook row + duplicate header leftovers) ---
$ws.Range("A22:O23").ClearContents()

$ws.Range("J21").Select()
